{"js": "// Replace the 25 division-expression placeholders with their updated values.\n// Every \"old\" string is unique in the document and every \"new\" string is\n// also unique, so a simple search-and-replace per pair is safe (no\n// re-matching of already-replaced text, no collisions between pairs).\nconst replacements = [\n  [\"191\u00f74=\", \"837\u00f79=\"],\n  [\"749\u00f78=\", \"768\u00f75=\"],\n  [\"226\u00f74=\", \"140\u00f73=\"],\n  [\"379\u00f73=\", \"878\u00f77=\"],\n  [\"118\u00f72=\", \"292\u00f76=\"],\n  [\"374\u00f75=\", \"634\u00f72=\"],\n  [\"797\u00f77=\", \"172\u00f72=\"],\n  [\"727\u00f76=\", \"767\u00f75=\"],\n  [\"385\u00f77=\", \"347\u00f77=\"],\n  [\"444\u00f78=\", \"833\u00f74=\"],\n  [\"847\u00f79=\", \"383\u00f78=\"],\n  [\"743\u00f73=\", \"326\u00f78=\"],\n  [\"718\u00f72=\", \"519\u00f79=\"],\n  [\"581\u00f79=\", \"429\u00f77=\"],\n  [\"345\u00f76=\", \"971\u00f77=\"],\n  [\"446\u00f79=\", \"540\u00f72=\"],\n  [\"985\u00f75=\", \"480\u00f72=\"],\n  [\"400\u00f76=\", \"209\u00f72=\"],\n  [\"427\u00f75=\", \"987\u00f76=\"],\n  [\"987\u00f77=\", \"479\u00f76=\"],\n  [\"519\u00f75=\", \"438\u00f73=\"],\n  [\"223\u00f73=\", \"499\u00f79=\"],\n  [\"755\u00f73=\", \"322\u00f79=\"],\n  [\"439\u00f75=\", \"907\u00f76=\"],\n  [\"688\u00f73=\", \"254\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-expression placeholders with their updated values.\n# Every \"old\" string is unique in the document and every \"new\" string is\n# also unique, so a simple Find/Replace-All per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"191\u00f74=\", \"837\u00f79=\"),\n    @(\"749\u00f78=\", \"768\u00f75=\"),\n    @(\"226\u00f74=\", \"140\u00f73=\"),\n    @(\"379\u00f73=\", \"878\u00f77=\"),\n    @(\"118\u00f72=\", \"292\u00f76=\"),\n    @(\"374\u00f75=\", \"634\u00f72=\"),\n    @(\"797\u00f77=\", \"172\u00f72=\"),\n    @(\"727\u00f76=\", \"767\u00f75=\"),\n    @(\"385\u00f77=\", \"347\u00f77=\"),\n    @(\"444\u00f78=\", \"833\u00f74=\"),\n    @(\"847\u00f79=\", \"383\u00f78=\"),\n    @(\"743\u00f73=\", \"326\u00f78=\"),\n    @(\"718\u00f72=\", \"519\u00f79=\"),\n    @(\"581\u00f79=\", \"429\u00f77=\"),\n    @(\"345\u00f76=\", \"971\u00f77=\"),\n    @(\"446\u00f79=\", \"540\u00f72=\"),\n    @(\"985\u00f75=\", \"480\u00f72=\"),\n    @(\"400\u00f76=\", \"209\u00f72=\"),\n    @(\"427\u00f75=\", \"987\u00f76=\"),\n    @(\"987\u00f77=\", \"479\u00f76=\"),\n    @(\"519\u00f75=\", \"438\u00f73=\"),\n    @(\"223\u00f73=\", \"499\u00f79=\"),\n    @(\"755\u00f73=\", \"322\u00f79=\"),\n    @(\"439\u00f75=\", \"907\u00f76=\"),\n    @(\"688\u00f73=\", \"254\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
